$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 8965.26
$ws.Range("B13").Value = 9011.2199999999993
$ws.Range("C13").Value = 17.8
$ws.Range("D13").Value = 17.89
$ws.Range("E13").Value = $true
$ws.Range("F13").Value = 0.51
$ws.Range("G13").Value = 42620.76630787037
$ws.Range("G13").NumberFormat = "m/d/yy h:mm"
$ws.Range("H13").Value = $false
